# Update the per-line power-flow results for the 380 kV case
# (Case_5_134, res_line/pl_mw). Columns B:F, I, K, L, N, O change on
# rows 2-25; columns A, G, H, J, M are untouched (index / zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.538671267783343
$ws.Range("C2").Value = 0.1051590893890193
$ws.Range("D2").Value = 0.041248916846925
$ws.Range("E2").Value = 0.09797392431060459
$ws.Range("F2").Value = 0.8658998088535341
$ws.Range("I2").Value = 0.7959050165789847
$ws.Range("K2").Value = 0.3331153853717694
$ws.Range("L2").Value = 0.2056283115882991
$ws.Range("N2").Value = 1.617569122752709
$ws.Range("O2").Value = 3.058396782020139

# Row 3
$ws.Range("B3").Value = 0.4977341740917325
$ws.Range("C3").Value = 0.1035833427650914
$ws.Range("D3").Value = 0.03904359781390099
$ws.Range("E3").Value = 0.09743708183723854
$ws.Range("F3").Value = 0.8657033429486631
$ws.Range("I3").Value = 0.8011007215043442
$ws.Range("K3").Value = 0.2959267575986644
$ws.Range("L3").Value = 0.1982594525952237
$ws.Range("N3").Value = 1.633924489781492
$ws.Range("O3").Value = 3.070715944059117

# Row 4
$ws.Range("B4").Value = 0.4727407509762145
$ws.Range("C4").Value = 0.1026049566053544
$ws.Range("D4").Value = 0.03767506608321725
$ws.Range("E4").Value = 0.09715595846720859
$ws.Range("F4").Value = 0.8660603216340874
$ws.Range("I4").Value = 0.8046630018676417
$ws.Range("K4").Value = 0.2731019816423839
$ws.Range("L4").Value = 0.1938435960834113
$ws.Range("N4").Value = 1.644483976683109
$ws.Range("O4").Value = 3.079838884716693

# Row 5
$ws.Range("B5").Value = 0.4625921552024579
$ws.Range("C5").Value = 0.1022035415682723
$ws.Range("D5").Value = 0.03711376696979585
$ws.Range("E5").Value = 0.09705361697469428
$ws.Range("F5").Value = 0.8663259775372723
$ws.Range("I5").Value = 0.8062082719799015
$ws.Range("K5").Value = 0.2638035503759824
$ws.Range("L5").Value = 0.1920715128021016
$ws.Range("N5").Value = 1.648917129387462
$ws.Range("O5").Value = 3.083948746751886

# Row 6
$ws.Range("B6").Value = 0.4609092090261697
$ws.Range("C6").Value = 0.1021367235026247
$ws.Range("D6").Value = 0.03702034626972051
$ws.Range("E6").Value = 0.09703736188988543
$ws.Range("F6").Value = 0.8663773508396631
$ws.Range("I6").Value = 0.8064705181707481
$ws.Range("K6").Value = 0.2622597409564662
$ws.Range("L6").Value = 0.1917789181870972
$ws.Range("N6").Value = 1.649661107784386
$ws.Range("O6").Value = 3.084654878968777

# Row 7
$ws.Range("B7").Value = 0.4726037350579873
$ws.Range("C7").Value = 0.1025995539488562
$ws.Range("D7").Value = 0.03766751079834307
$ws.Range("E7").Value = 0.09715452875119723
$ws.Range("F7").Value = 0.8660634176284034
$ws.Range("I7").Value = 0.8046834628297468
$ws.Range("K7").Value = 0.2729765674920088
$ws.Range("L7").Value = 0.193819586026919
$ws.Range("N7").Value = 1.644543237036654
$ws.Range("O7").Value = 3.079892723548213

# Row 8
$ws.Range("B8").Value = 0.5245270623254044
$ws.Range("C8").Value = 0.1046180392462261
$ws.Range("D8").Value = 0.04049153719945764
$ws.Range("E8").Value = 0.09777876813330977
$ws.Range("F8").Value = 0.8657329883880962
$ws.Range("I8").Value = 0.7976192937273936
$ws.Range("K8").Value = 0.3202911507573276
$ws.Range("L8").Value = 0.2030650309879718
$ws.Range("N8").Value = 1.62310111396725
$ws.Range("O8").Value = 3.062320980859653

# Row 9
$ws.Range("B9").Value = 0.6274510899459642
$ws.Range("C9").Value = 0.1084893561278051
$ws.Range("D9").Value = 0.04591390937508777
$ws.Range("E9").Value = 0.09938689841980519
$ws.Range("F9").Value = 0.8688720300209738
$ws.Range("I9").Value = 0.7867171728087534
$ws.Range("K9").Value = 0.4131286422375524
$ws.Range("L9").Value = 0.2220544979287951
$ws.Range("N9").Value = 1.585155804068794
$ws.Range("O9").Value = 3.04022632527213

# Row 10
$ws.Range("B10").Value = 0.703715481109441
$ws.Range("C10").Value = 0.1112799672328038
$ws.Range("D10").Value = 0.04982651940309069
$ws.Range("E10").Value = 0.1008016338701445
$ws.Range("F10").Value = 0.8734851428682902
$ws.Range("I10").Value = 0.7805042750361579
$ws.Range("K10").Value = 0.4813499111045303
$ws.Range("L10").Value = 0.2365275647398022
$ws.Range("N10").Value = 1.559775730341798
$ws.Range("O10").Value = 3.031525887103356

# Row 11
$ws.Range("B11").Value = 0.7385453508412638
$ws.Range("C11").Value = 0.1125377160006309
$ws.Range("D11").Value = 0.05159085896970339
$ws.Range("E11").Value = 0.101495714173776
$ws.Range("F11").Value = 0.8760844041170444
$ws.Range("I11").Value = 0.7780676722788868
$ws.Range("K11").Value = 0.5123847872752378
$ws.Range("L11").Value = 0.2432245889232263
$ws.Range("N11").Value = 1.548771592419298
$ws.Range("O11").Value = 3.029202715278331

# Row 12
$ws.Range("B12").Value = 0.7517535861373119
$ws.Range("C12").Value = 0.1130122920402101
$ws.Range("D12").Value = 0.05225671511964691
$ws.Range("E12").Value = 0.1017657882381862
$ws.Range("F12").Value = 0.8771406247382245
$ws.Range("I12").Value = 0.7772009957414951
$ws.Range("K12").Value = 0.5241365198428412
$ws.Range("L12").Value = 0.2457767771810069
$ws.Range("N12").Value = 1.544682453440084
$ws.Range("O12").Value = 3.028557945559299

# Row 13
$ws.Range("B13").Value = 0.7489081269086455
$ws.Range("C13").Value = 0.1129101598233149
$ws.Range("D13").Value = 0.05211341209275133
$ws.Range("E13").Value = 0.1017073012718903
$ws.Range("F13").Value = 0.8769099504141735
$ws.Range("I13").Value = 0.7773851592539955
$ws.Range("K13").Value = 0.5216056058189906
$ws.Range("L13").Value = 0.2452264001990301
$ws.Range("N13").Value = 1.545559657896668
$ws.Range("O13").Value = 3.028686359773985

# Row 14
$ws.Range("B14").Value = 0.739631625456127
$ws.Range("C14").Value = 0.1125767939738296
$ws.Range("D14").Value = 0.05164568483124299
$ws.Range("E14").Value = 0.101517788339379
$ws.Range("F14").Value = 0.8761698586956825
$ws.Range("I14").Value = 0.7779952478446788
$ws.Range("K14").Value = 0.5133516228712836
$ws.Range("L14").Value = 0.2434342354854522
$ws.Range("N14").Value = 1.548433615247763
$ws.Range("O14").Value = 3.02914496123023

# Row 15
$ws.Range("B15").Value = 0.7339519419436726
$ws.Range("C15").Value = 0.1123723749368821
$ws.Range("D15").Value = 0.05135889309521247
$ws.Range("E15").Value = 0.1014026485114243
$ws.Range("F15").Value = 0.8757258974080528
$ws.Range("I15").Value = 0.7783762384658388
$ws.Range("K15").Value = 0.5082957368848895
$ws.Range("L15").Value = 0.2423385857288309
$ws.Range("N15").Value = 1.55020414177678
$ws.Range("O15").Value = 3.029456464170408

# Row 16
$ws.Range("B16").Value = 0.7014419558196323
$ws.Range("C16").Value = 0.111197532824427
$ws.Range("D16").Value = 0.04971090081727425
$ws.Range("E16").Value = 0.1007572885720478
$ws.Range("F16").Value = 0.8733253452706435
$ws.Range("I16").Value = 0.7806713520752169
$ws.Range("K16").Value = 0.479321671420962
$ws.Range("L16").Value = 0.2360921667980875
$ws.Range("N16").Value = 1.560505773939852
$ws.Range("O16").Value = 3.031710597432806

# Row 17
$ws.Range("B17").Value = 0.6815326375670452
$ws.Range("C17").Value = 0.1104737895494665
$ws.Range("D17").Value = 0.04869591416915853
$ws.Range("E17").Value = 0.1003743021107333
$ws.Range("F17").Value = 0.8719808859547058
$ws.Range("I17").Value = 0.782179118971797
$ws.Range("K17").Value = 0.4615467685254373
$ws.Range("L17").Value = 0.2322891009502257
$ws.Range("N17").Value = 1.56696418006926
$ws.Range("O17").Value = 3.03351204842636

# Row 18
$ws.Range("B18").Value = 0.6700942471657072
$ws.Range("C18").Value = 0.1100564102530228
$ws.Range("D18").Value = 0.04811066097735051
$ws.Range("E18").Value = 0.1001587740179311
$ws.Range("F18").Value = 0.8712547240074144
$ws.Range("I18").Value = 0.7830830260889812
$ws.Range("K18").Value = 0.451323214253307
$ws.Range("L18").Value = 0.2301123346785232
$ws.Range("N18").Value = 1.570729839693509
$ws.Range("O18").Value = 3.034702067418635

# Row 19
$ws.Range("B19").Value = 0.6662236472070049
$ws.Range("C19").Value = 0.1099149043035013
$ws.Range("D19").Value = 0.04791225445401182
$ws.Range("E19").Value = 0.1000866173657116
$ws.Range("F19").Value = 0.8710169556729639
$ws.Range("I19").Value = 0.7833953739393067
$ws.Range("K19").Value = 0.4478617298171059
$ws.Range("L19").Value = 0.2293771523636963
$ws.Range("N19").Value = 1.572013577811679
$ws.Range("O19").Value = 3.035131418196073

# Row 20
$ws.Range("B20").Value = 0.6836506864756302
$ws.Range("C20").Value = 0.1105509473938824
$ws.Range("D20").Value = 0.04880411251070882
$ws.Range("E20").Value = 0.1004145796493567
$ws.Range("F20").Value = 0.8721191280869647
$ws.Range("I20").Value = 0.7820148185882942
$ws.Range("K20").Value = 0.4634389312974463
$ws.Range("L20").Value = 0.2326928416105289
$ws.Range("N20").Value = 1.566271397761112
$ws.Range("O20").Value = 3.033304356455943

# Row 21
$ws.Range("B21").Value = 0.7423558513435751
$ws.Range("C21").Value = 0.1126747580817096
$ws.Range("D21").Value = 0.05178312925298201
$ws.Range("E21").Value = 0.1015732565439542
$ws.Range("F21").Value = 0.8763852897762305
$ws.Range("I21").Value = 0.7778145300681487
$ws.Range("K21").Value = 0.5157760361308306
$ws.Range("L21").Value = 0.2439601996475886
$ws.Range("N21").Value = 1.547587350095868
$ws.Range("O21").Value = 3.029003882809121

# Row 22
$ws.Range("B22").Value = 0.7808329971532544
$ws.Range("C22").Value = 0.1140528365433795
$ws.Range("D22").Value = 0.0537168971473676
$ws.Range("E22").Value = 0.1023727140169903
$ws.Range("F22").Value = 0.8795927435221103
$ws.Range("I22").Value = 0.7753958691518505
$ws.Range("K22").Value = 0.5499781580378738
$ws.Range("L22").Value = 0.2514182764704174
$ws.Range("N22").Value = 1.535830262570573
$ws.Range("O22").Value = 3.027562782080707

# Row 23
$ws.Range("B23").Value = 0.7602872114176478
$ws.Range("C23").Value = 0.1133182481173876
$ws.Range("D23").Value = 0.05268602524598975
$ws.Range("E23").Value = 0.101942175080076
$ws.Range("F23").Value = 0.8778425245390622
$ws.Range("I23").Value = 0.7766568894895798
$ws.Range("K23").Value = 0.5317243329693042
$ws.Range("L23").Value = 0.2474291740566628
$ws.Range("N23").Value = 1.542063685384085
$ws.Range("O23").Value = 3.028206648975782

# Row 24
$ws.Range("B24").Value = 0.6826930923919008
$ws.Range("C24").Value = 0.1105160683486588
$ws.Range("D24").Value = 0.04875520140789291
$ws.Range("E24").Value = 0.1003963556700604
$ws.Range("F24").Value = 0.8720564830803781
$ws.Range("I24").Value = 0.7820889833135567
$ws.Range("K24").Value = 0.46258349856069
$ws.Range("L24").Value = 0.2325102803402643
$ws.Range("N24").Value = 1.56658444068636
$ws.Range("O24").Value = 3.033397773225801

# Row 25
$ws.Range("B25").Value = 0.5994920678619451
$ws.Range("C25").Value = 0.1074514360040411
$ws.Range("D25").Value = 0.04445946247155774
$ws.Range("E25").Value = 0.0989108347578167
$ws.Range("F25").Value = 0.8676175942544475
$ws.Range("I25").Value = 0.7893507933879711
$ws.Range("K25").Value = 0.3880098704332511
$ws.Range("L25").Value = 0.2168255955495653
$ws.Range("N25").Value = 1.594982092841879
$ws.Range("O25").Value = 3.044880393637953
